$wb = $excel.ActiveWorkbook

# --- Sheet 1: LP1912 ---
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Range("A2").Value = "Última actualización: 02:56:06"
$ws1.Range("A3").Value = "Total filas: 10"

$ws1.Cells.Item(14, 1).Value = "02:56:05"
$ws1.Cells.Item(14, 2).Value = "04:46"
$ws1.Cells.Item(14, 3).Value = "215A_EL PATO"
$ws1.Cells.Item(14, 4).Value = 110
$ws1.Cells.Item(14, 5).Value = "LP1912"

$ws1.Cells.Item(15, 1).Value = "02:56:05"
$ws1.Cells.Item(15, 2).Value = "04:53"
$ws1.Cells.Item(15, 3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(15, 4).Value = 117
$ws1.Cells.Item(15, 5).Value = "LP1912"

# --- Sheet 2: LP1912-215 ---
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: 02:56:06"
$ws2.Range("A3").Value = "Total filas: 4"

$ws2.Cells.Item(9, 1).Value = "02:56:05"
$ws2.Cells.Item(9, 2).Value = "04:46"
$ws2.Cells.Item(9, 3).Value = "215A_EL PATO"
$ws2.Cells.Item(9, 4).Value = 110
$ws2.Cells.Item(9, 5).Value = "LP1912"

# --- Sheet 3: 6203-6173 ---
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: 02:56:06"
